# FEAT #7647 change to contact.postal_address in templates
$d = $word.ActiveDocument

# --- 1) Table cell margins: left padding 108 dxa (5.4pt) -> 113 dxa (5.65pt) on both tables ---
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $t.LeftPadding = 5.65
}

# --- 2) Placeholder text: [contact.afnor] -> [contact.postal_address] ---
$d.Content.Find.Execute("[contact.afnor]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[contact.postal_address]", 2)

# --- 3) Cached TIME field result: 17/04/2018 -> 25/04/2018 ---
$d.Content.Find.Execute("17/04/2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25/04/2018", 2)

# --- 4) Header horizontal-line shape: give it a name ("" -> "Image1") and nudge its
#        size very slightly (the same no-op-looking resize Word performs when it
#        re-touches the shape), which is how the underlying extent values move. ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$lineShape = $hdr.Shapes.Item(2)
$lineShape.Name = "Image1"
$lineShape.Width = 543.7
$lineShape.Height = 0.15

# --- 5) Style defaults: Normal style font color auto -> 00000A (RRGGBB -> VBA BGR long) ---
$normal = $d.Styles.Item("Normal")
$normal.Font.Color = 655360
